$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 (Marking): Right count 4 -> 5, Wrong penalty -1 -> -1.2
$ws.Range("B11").Value = 5
$ws.Range("C11").Value = -1.2

# Row 12 (Total): Right total 16 -> 20, Wrong total 0 -> -0, fraction text update
$ws.Range("B12").Value = 20
$ws.Range("C12").Value = -0
$ws.Range("E12").Value = "20.0/140"
